$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Sheet1 (ActivityLog): fill in two previously-empty cells in the Sat row (row 6)
# ---------------------------------------------------------------------------
$wsLog = $wb.Worksheets.Item("ActivityLog")
$wsLog.Range("F6").Value = "Contribute to repository via Git"
$wsLog.Range("H6").Value = "Learn Google App Engine code"

# ---------------------------------------------------------------------------
# Sheet2 (ProductBacklog): build the product backlog table
# String cell values are written in the exact order needed so that the
# resulting shared-strings table is ordered the same way as the target file.
# ---------------------------------------------------------------------------
$wsPb = $wb.Worksheets.Item("ProductBacklog")

$wsPb.Range("C2").Value = "ProuctBacklog"
$wsPb.Range("D2").Value = "User Story"
$wsPb.Range("E2").Value = "Sprint"
$wsPb.Range("C3").Value = "Current Location"
$wsPb.Range("B2").Value = "Pb_id"
$wsPb.Range("D3").Value = "I want to see my current location on google map."
$wsPb.Range("C6").Value = "Taxi Location"
$wsPb.Range("D6").Value = "I want to see nearby available taxis on google map"
$wsPb.Range("F2").Value = "Delivery Date"

$wsPb.Range("B3").Value = 1
$wsPb.Range("E3").Value = 1
$wsPb.Range("E4").Value = 1
$wsPb.Range("E5").Value = 1
$wsPb.Range("E6").Value = 1
$wsPb.Range("E7").Value = 2
$wsPb.Range("E8").Value = 2
$wsPb.Range("E9").Value = 2
$wsPb.Range("E10").Value = 2

$wsPb.Columns("C:C").ColumnWidth = 15.7109375
$wsPb.Columns("D:D").ColumnWidth = 45
$wsPb.Columns("F:F").ColumnWidth = 13.140625

$wsPb.Range("F3").Select() | Out-Null

# ---------------------------------------------------------------------------
# Sheet3 (SprintBacklog): build the sprint backlog table
# ---------------------------------------------------------------------------
$wsSb = $wb.Worksheets.Item("SprintBacklog")

$wsSb.Range("B2").Value = "Sb_id"
$wsSb.Range("D2").Value = "SprintBacklog"
$wsSb.Range("D7").Value = "TaxiAvailability API header"
$wsSb.Range("E2").Value = "Description"
$wsSb.Range("E7").Value = "Need to get API key and other header information"
$wsSb.Range("D3").Value = "JINJA2"
$wsSb.Range("E3").Value = "Use JINJA2 template for MVC framework"
$wsSb.Range("D4").Value = "JavaScript"
$wsSb.Range("E4").Value = "Implement current location in JavaScript in accordance with Google Map API"

$wsSb.Range("C2").Value = "Pb_id"
$wsSb.Range("F2").Value = "Sprint"
$wsSb.Range("G2").Value = "Delivery Date"

$wsSb.Range("B3").Value = 1
$wsSb.Range("C3").Value = 1
$wsSb.Range("F3").Value = 1
$wsSb.Range("B4").Value = 2
$wsSb.Range("C4").Value = 1
$wsSb.Range("F4").Value = 1
$wsSb.Range("B5").Value = 3
$wsSb.Range("C5").Value = 1
$wsSb.Range("F5").Value = 1
$wsSb.Range("C6").Value = 4

$wsSb.Columns("B:B").ColumnWidth = 5.85546875
$wsSb.Columns("C:C").ColumnWidth = 6
$wsSb.Columns("D:D").ColumnWidth = 25.140625
$wsSb.Columns("E:E").ColumnWidth = 69.85546875
$wsSb.Columns("F:F").ColumnWidth = 6.28515625
$wsSb.Columns("G:G").ColumnWidth = 13.140625

$wsSb.Range("D2").Select() | Out-Null

# ---------------------------------------------------------------------------
# Sheet1: move the "Sprint2" label textbox down, and refresh the selection
# ---------------------------------------------------------------------------
$wsLog.Activate()
$shp = $wsLog.Shapes.Item("TextBox 4")
$shp.Top = 1052.05

$wsLog.Range("H6").Select() | Out-Null
